$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Execution_Sheet")
$ws2 = $wb.Worksheets.Item("Credentials")

# --- Execution_Sheet: toggle YES/NO flags in column A ---
$ws1.Range("A2").Value = "NO"
$ws1.Range("A4").Value = "NO"
$ws1.Range("A7").Value = "NO"
$ws1.Range("A8").Value = "YES"
$ws1.Range("A13").Value = "NO"
$ws1.Range("A15").Value = "NO"
$ws1.Range("A17").Value = "NO"
$ws1.Range("A19").Value = "NO"
$ws1.Range("A21").Value = "NO"
$ws1.Range("A22").Value = "NO"
$ws1.Range("A24").Value = "NO"
$ws1.Range("A25").Value = "NO"
$ws1.Range("A27").Value = "NO"

# Rows 28-39: value becomes "NO" and the cell format in column A is
# normalized to match the same style used by the other A-column cells
# (copy format from A2, which already carries that style).
$ws1.Range("A28:A39").Value = "NO"
$ws1.Range("A2").Copy()
$ws1.Range("A28:A39").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Execution_Sheet: dual SIM / SMS handling updates (rows 8-12) ---
# Relabel the rework tags on rows 10-12 first ...
$ws1.Range("B10").Value = "Rework_60"
$ws1.Range("B11").Value = "Rework_61"
$ws1.Range("B12").Value = "Rework_62"

# ... then switch the SMS scenario rows to use the dual-SIM MSISDN (9732)
# and its corresponding "more data" shared-string labels.
$ws1.Range("D8").Value = "LIVE_USAGE_SMS"
$ws1.Range("K8").Value = "9732"
$ws1.Range("L8").Value = "971520001714 Moredata20"

$ws1.Range("D9").Value = "LIVE_USAGE_SMS"
$ws1.Range("K9").Value = "9732"
$ws1.Range("L9").Value = "971520001714 Moredata50"

$ws1.Range("D10").Value = "LIVE_USAGE_SMS"
$ws1.Range("K10").Value = "9732"
$ws1.Range("L10").Value = "971520001714 Moredata100"

$ws1.Range("D11").Value = "LIVE_USAGE_SMS"
$ws1.Range("K11").Value = "9732"
$ws1.Range("L11").Value = "971520001714 Moredata200"

$ws1.Range("D12").Value = "LIVE_USAGE_SMS"
$ws1.Range("K12").Value = "9732"
$ws1.Range("L12").Value = "971520001714 Moredata500"

# --- Credentials sheet: update F2:F8 values from 0 to 30 ---
$ws2.Range("F2:F8").Value = 30

# --- Selections ---
$ws1.Range("A8").Select()
$ws2.Range("F2:F8").Select()
$ws1.Activate()
